$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("2023-10-07 20:12:06", "hatespeech", "setting1", 30, 0.2181525160943794)
    ,@("2023-10-07 20:12:08", "hatespeech", "def", 40, 0.1770474220080296)
    ,@("2023-10-07 20:12:12", "hatespeech", "setting2", 10, 0.3687429832204173)
    ,@("2023-10-07 20:12:17", "hatespeech", "def", 30, 0.216860808076717)
    ,@("2023-10-07 20:12:17", "hatespeech", "setting1", 10, 0.3663133674668461)
    ,@("2023-10-07 20:12:19", "hatespeech", "def", 10, 0.3663292266055626)
    ,@("2023-10-07 20:12:20", "hatespeech", "setting6", 40, 0.1816254702245944)
    ,@("2023-10-07 20:12:20", "hatespeech", "setting3", 10, 0.3647653896434879)
    ,@("2023-10-07 20:12:20", "hatespeech", "setting2", 20, 0.2776621958449362)
    ,@("2023-10-07 20:12:21", "hatespeech", "def", 20, 0.2748046669632554)
    ,@("2023-10-07 20:12:21", "hatespeech", "setting4", 20, 0.278008186185377)
    ,@("2023-10-07 20:12:21", "hatespeech", "setting2", 30, 0.2207366235281553)
    ,@("2023-10-07 20:12:22", "hatespeech", "setting5", 40, 0.1825340182958585)
    ,@("2023-10-07 20:12:22", "hatespeech", "setting4", 40, 0.1862651554593032)
    ,@("2023-10-07 20:12:22", "hatespeech", "setting2", 40, 0.1781252847485837)
    ,@("2023-10-07 20:12:24", "hatespeech", "setting5", 10, 0.368563903575915)
    ,@("2023-10-07 20:12:24", "hatespeech", "setting6", 10, 0.366636673618901)
    ,@("2023-10-07 20:12:25", "hatespeech", "setting5", 30, 0.2264740281900377)
    ,@("2023-10-07 20:12:26", "hatespeech", "setting4", 10, 0.3676986503657903)
    ,@("2023-10-07 20:12:27", "hatespeech", "setting3", 20, 0.2711620338916789)
    ,@("2023-10-07 20:12:27", "hatespeech", "setting1", 40, 0.1748688091151119)
    ,@("2023-10-07 20:12:27", "hatespeech", "setting3", 40, 0.1781794415245524)
    ,@("2023-10-07 20:12:27", "hatespeech", "setting3", 30, 0.2204962227375961)
    ,@("2023-10-07 20:12:27", "hatespeech", "setting5", 20, 0.2782320654789707)
    ,@("2023-10-07 20:12:28", "hatespeech", "setting6", 20, 0.2748596905958181)
    ,@("2023-10-07 20:12:28", "hatespeech", "setting1", 20, 0.2717033399253748)
    ,@("2023-10-07 20:12:28", "hatespeech", "setting6", 30, 0.2222684071561982)
    ,@("2023-10-07 20:12:28", "hatespeech", "setting4", 30, 0.2230270797175439)
)

$startRow = 206
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
}
